$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "Hostius"
$ws.Range("B11").Select()
